$d = $word.ActiveDocument

$pairs = @(
    @("11×63=", "23×51="),
    @("73×88=", "11×93="),
    @("58×41=", "25×67="),
    @("52×14=", "53×86="),
    @("37×37=", "15×32="),
    @("17×29=", "13×26="),
    @("66×68=", "72×46="),
    @("69×93=", "71×73="),
    @("79×92=", "82×91="),
    @("58×17=", "76×50="),
    @("41×44=", "46×61="),
    @("91×69=", "19×80="),
    @("28×93=", "34×38="),
    @("40×72=", "34×32="),
    @("36×54=", "45×21="),
    @("75×66=", "12×97="),
    @("36×62=", "83×67="),
    @("98×16=", "97×50="),
    @("81×54=", "95×76="),
    @("20×61=", "42×55="),
    @("88×94=", "15×21="),
    @("67×52=", "54×58="),
    @("97×12=", "94×80="),
    @("38×52=", "81×13="),
    @("79×78=", "55×98=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
